$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172157764434814
$ws.Range("B1").Value = 2.385447978973389
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.371347188949585
$ws.Range("E1").Value = 1.209865927696228
